$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, shifting existing rows 55-165 down to 56-166
$ws.Rows(55).Insert()

# Populate the newly inserted row 55 with the new weekly data point
$ws.Cells.Item(55, 1).Value = 7
$ws.Cells.Item(55, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(55, 3).Value = "Ñuble"
$ws.Cells.Item(55, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = 16
$ws.Cells.Item(55, 6).Value = 100112003
$ws.Cells.Item(55, 7).Value = "Ajo"
$ws.Cells.Item(55, 8).Value = "Chino"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 60
$ws.Cells.Item(55, 11).Value = 18000
$ws.Cells.Item(55, 12).Value = 19000
$ws.Cells.Item(55, 13).Value = 18500
$ws.Cells.Item(55, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(55, 15).Value = "China"
$ws.Cells.Item(55, 16).Value = 1850
$ws.Cells.Item(55, 17).Value = 10
$ws.Cells.Item(55, 18).Value = "Hortaliza"
